$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 742.8421
$ws.Range("I80").Value = 628.5454999999999
$ws.Range("J80").Value = 900
$ws.Range("K80").Value = 1885.6365
$ws.Range("L80").Value = 2700
$ws.Range("M80").Value = -887.6364999999998
$ws.Range("N80").Value = -4696

$ws.Range("H83").Value = 742.8421
$ws.Range("I83").Value = 628.5454999999999
$ws.Range("J83").Value = 900
$ws.Range("K83").Value = 5656.9095
$ws.Range("L83").Value = 8100
$ws.Range("M83").Value = -664.9094999999998
$ws.Range("N83").Value = -18084

$ws.Range("H86").Value = 3613.35
$ws.Range("I86").Value = 4393.3076
$ws.Range("K86").Value = 4393.3076
$ws.Range("M86").Value = -3270.3076

$ws.Range("H89").Value = 3613.35
$ws.Range("I89").Value = 4393.3076
$ws.Range("K89").Value = 21966.538
$ws.Range("M89").Value = -16350.538

$ws.Range("H94").Value = 2400
$ws.Range("I94").Value = 2400
$ws.Range("K94").Value = 2400
$ws.Range("M94").Value = -1949

$ws.Range("H107").Value = 25006432
$ws.Range("I107").Value = 35715772
$ws.Range("J107").Value = 17966.666
$ws.Range("K107").Value = 35715772
$ws.Range("L107").Value = 17966.666
$ws.Range("M107").Value = -35713852
$ws.Range("N107").Value = -21806.666

$ws.Range("H141").Value = 2721.889
$ws.Range("I141").Value = 2721.889
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 8165.667
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -2985.667
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 21142.4
$ws.Range("I45").Value = 21142.4
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 21142.4
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -20765.4
$ws.Range("N45").ClearContents()

$ws.Range("H74").Value = 1853.8182
$ws.Range("I74").Value = 1792.7742
$ws.Range("K74").Value = 1792.7742
$ws.Range("M74").Value = -918.7742000000001

$ws.Range("H77").Value = 1853.8182
$ws.Range("I77").Value = 1792.7742
$ws.Range("K77").Value = 8963.871000000001
$ws.Range("M77").Value = -4595.871000000001

$ws.Range("H123").Value = 45429
$ws.Range("J123").Value = 45429
$ws.Range("L123").Value = 45429
$ws.Range("N123").Value = -55229

$ws.Range("H140").Value = 38883.332
$ws.Range("J140").Value = 38883.332
$ws.Range("L140").Value = 38883.332
$ws.Range("N140").Value = -49243.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 126.64286
$ws.Range("I7").Value = 118.625
$ws.Range("J7").Value = 137.33333
$ws.Range("K7").Value = 118.625
$ws.Range("L7").Value = 137.33333
$ws.Range("M7").Value = -5.625
$ws.Range("N7").Value = -363.33333

$ws.Range("H122").Value = 903.8570999999999
$ws.Range("I122").Value = 935.5
$ws.Range("J122").Value = 714
$ws.Range("K122").Value = 2806.5
$ws.Range("L122").Value = 2142
$ws.Range("M122").Value = -356.5
$ws.Range("N122").Value = -7042

$ws.Range("H134").Value = 3338.625
$ws.Range("I134").Value = 3710.4736
$ws.Range("K134").Value = 11131.4208
$ws.Range("M134").Value = -8596.4208

$ws.Range("H135").Value = 31927.646
$ws.Range("J135").Value = 31927.646
$ws.Range("L135").Value = 31927.646
$ws.Range("N135").Value = -42067.646

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 670.75
$ws.Range("I2").Value = 925.4545000000001
$ws.Range("J2").Value = 110.4
$ws.Range("K2").Value = 5552.727000000001
$ws.Range("L2").Value = 662.4000000000001
$ws.Range("M2").Value = -5439.727000000001
$ws.Range("N2").Value = -888.4000000000001

$ws.Range("H42").Value = 5000
$ws.Range("J42").Value = 5000
$ws.Range("L42").Value = 15000
$ws.Range("N42").Value = -16068

$ws.Range("H54").Value = 2958.3333
$ws.Range("J54").Value = 2958.3333
$ws.Range("L54").Value = 8874.999899999999
$ws.Range("N54").Value = -9992.999899999999

$ws.Range("H113").Value = 1936043.6
$ws.Range("I113").Value = 2941665.2
$ws.Range("J113").Value = 714931.7
$ws.Range("K113").Value = 8824995.600000001
$ws.Range("L113").Value = 2144795.1
$ws.Range("M113").Value = -8822825.600000001
$ws.Range("N113").Value = -2149135.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 33691
$ws.Range("J15").Value = 33691
$ws.Range("L15").Value = 33691
$ws.Range("N15").Value = -34267

$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

$ws.Range("H81").Value = 33691
$ws.Range("J81").Value = 33691
$ws.Range("L81").Value = 33691
$ws.Range("N81").Value = -35687

$ws.Range("H84").Value = 33691
$ws.Range("J84").Value = 33691
$ws.Range("L84").Value = 101073
$ws.Range("N84").Value = -111057

$ws.Range("H122").Value = 5742616
$ws.Range("I122").Value = 5893892.5
$ws.Range("J122").Value = 5557722
$ws.Range("K122").Value = 17681677.5
$ws.Range("L122").Value = 16673166
$ws.Range("M122").Value = -17679227.5
$ws.Range("N122").Value = -16678066

$ws.Range("H123").Value = 23027.23
$ws.Range("J123").Value = 23027.23
$ws.Range("L123").Value = 23027.23
$ws.Range("N123").Value = -27927.23

$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1299.2222
$ws.Range("I16").Value = 840.8
$ws.Range("J16").Value = 1872.25
$ws.Range("K16").Value = 840.8
$ws.Range("L16").Value = 1872.25
$ws.Range("M16").Value = -670.8
$ws.Range("N16").Value = -2212.25

$ws.Range("H22").Value = 2138352.2
$ws.Range("I22").Value = 11111366
$ws.Range("J22").Value = 1920.4524
$ws.Range("K22").Value = 11111366
$ws.Range("L22").Value = 1920.4524
$ws.Range("M22").Value = -11111071
$ws.Range("N22").Value = -2510.4524

$ws.Range("H27").Value = 2138352.2
$ws.Range("I27").Value = 11111366
$ws.Range("J27").Value = 1920.4524
$ws.Range("K27").Value = 11111366
$ws.Range("L27").Value = 1920.4524
$ws.Range("M27").Value = -11111259
$ws.Range("N27").Value = -2134.4524

$ws.Range("H40").Value = 90911050
$ws.Range("I40").Value = 125001816
$ws.Range("K40").Value = 125001816
$ws.Range("M40").Value = -125001680

$ws.Range("H68").Value = 71430750
$ws.Range("I68").Value = 1703
$ws.Range("K68").Value = 1703
$ws.Range("M68").Value = -954

$ws.Range("H71").Value = 71430750
$ws.Range("I71").Value = 1703
$ws.Range("K71").Value = 8515
$ws.Range("M71").Value = -4771

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 7558.2144
$ws.Range("I100").Value = 9465.272000000001
$ws.Range("J100").Value = 565.6667
$ws.Range("K100").Value = 18930.544
$ws.Range("L100").Value = 1131.3334
$ws.Range("M100").Value = -18389.544
$ws.Range("N100").Value = -2213.3334

